# Applies the "Updated cryptos list" data refresh to Sheet1.
# Column D (Price) and E (Volume(1h)) values are updated for most rows;
# rows 40-47 additionally have their Coin (B) and Link (C) values updated
# because the coin ranking order changed.
#
# Values in column D that look like plain numbers ("1.00", "403.67", ...)
# are written with a leading apostrophe so Excel keeps them as text
# (matching the source data, which stores everything as strings - e.g.
# "64.310.67" is not a valid number, and "1.00" must keep its trailing zero).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.008.29"
$ws.Range("E2").Value = "  -3.93%  "
$ws.Range("D3").Value = "3.611.71"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").Value = "'403.67"
$ws.Range("E5").Value = "  -2.83%  "
$ws.Range("D6").Value = "'131.36"
$ws.Range("E6").Value = "  +1.20%  "
$ws.Range("D7").Value = "3.606.69"
$ws.Range("E7").Value = "  +0.67%  "
$ws.Range("D8").Value = "'0.616"
$ws.Range("E8").Value = "  -5.14%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").Value = "'0.720"
$ws.Range("E10").Value = "  -7.14%  "
$ws.Range("D11").Value = "'0.158"
$ws.Range("E11").Value = "  -9.64%  "
$ws.Range("D12").Value = "'0.0000318"
$ws.Range("E12").Value = "  -4.00%  "
$ws.Range("D13").Value = "'41.44"
$ws.Range("E13").Value = "  -2.20%  "
$ws.Range("D14").Value = "'9.85"
$ws.Range("E14").Value = "  -0.30%  "
$ws.Range("D15").Value = "4.208.25"
$ws.Range("E15").Value = "  +1.81%  "
$ws.Range("E16").Value = "  -1.18%  "
$ws.Range("D17").Value = "3.606.37"
$ws.Range("E17").Value = "  +0.89%  "
$ws.Range("D18").Value = "'13.44"
$ws.Range("E18").Value = "  +9.34%  "
$ws.Range("D19").Value = "'19.83"
$ws.Range("E19").Value = "  -1.96%  "
$ws.Range("E20").Value = "  -4.93%  "
$ws.Range("D21").Value = "64.183.13"
$ws.Range("E21").Value = "  -3.58%  "
$ws.Range("D22").Value = "'418.33"
$ws.Range("E22").Value = "  -6.33%  "
$ws.Range("D23").Value = "'15.04"
$ws.Range("E23").Value = "  +15.23%  "
$ws.Range("D24").Value = "'85.06"
$ws.Range("E24").Value = "  -4.58%  "
$ws.Range("E25").Value = "  -5.57%  "
$ws.Range("D26").Value = "'35.41"
$ws.Range("E26").Value = "  +1.95%  "
$ws.Range("D27").Value = "'3.17"
$ws.Range("E27").Value = "  -4.46%  "
$ws.Range("D28").Value = "'9.32"
$ws.Range("E28").Value = "  -6.60%  "
$ws.Range("D29").Value = "'5.14"
$ws.Range("E29").Value = "  +5.87%  "
$ws.Range("D30").Value = "'12.75"
$ws.Range("E30").Value = "  +3.25%  "
$ws.Range("E31").Value = "  -2.23%  "
$ws.Range("E32").Value = "  -1.36%  "
$ws.Range("D33").Value = "'6.90"
$ws.Range("E33").Value = "  -6.03%  "
$ws.Range("D34").Value = "'40.96"
$ws.Range("E34").Value = "  +2.84%  "
$ws.Range("D35").Value = "'0.159"
$ws.Range("E35").Value = "  -0.20%  "
$ws.Range("D36").Value = "'55.55"
$ws.Range("E36").Value = "  -1.78%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").Value = "'0.0461"
$ws.Range("E38").Value = "  -5.91%  "
$ws.Range("D39").Value = "'2.89"
$ws.Range("E39").Value = "  +27.18%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").Value = "'0.998"
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").Value = "'0.139"
$ws.Range("E41").Value = "  -5.43%  "
$ws.Range("B42").Value = "ApeXProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D42").Value = "'3.14"
$ws.Range("E42").Value = "  +22.96%  "
$ws.Range("B43").Value = "NEARProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D43").Value = "'4.39"
$ws.Range("E43").Value = "  +1.98%  "
$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").Value = "0.0₃0634"
$ws.Range("E44").Value = "  -13.21%  "
$ws.Range("D45").Value = "'3.28"
$ws.Range("E45").Value = "  +1.16%  "
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").Value = "'144.40"
$ws.Range("E46").Value = "  -2.86%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'26.43"
$ws.Range("E47").Value = "  +25.35%  "
$ws.Range("E48").Value = "  +4.64%  "
$ws.Range("D49").Value = "'2.80"
$ws.Range("E49").Value = "  -5.96%  "
$ws.Range("E50").Value = "  -7.04%  "
$ws.Range("D51").Value = "'0.289"
$ws.Range("E51").Value = "  -6.34%  "
